# Applies the diff to the one-slide poster deck:
#  - reposition the "Time Series Analysis" header box (shape 102)
#  - reposition/resize its body textbox (shape 103)
#  - two text tweaks in the "Linear Regression" body box (shape 105)
#  - one text tweak in the "Conclusion" body box (shape 107)
#  - reposition/resize a picture (shape 108) and restyle the table (shape 109)
#  - reposition/resize two more pictures (shapes 110, 111)
#  - reposition the "Linear Regression" header box (shape 112)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 102 ("Time S eries Analysis" header) : move only ---
$sh = $s.Shapes.Item(7)
$sh.Left = 719.8641967773438
$sh.Top = 829.1437377929688

# --- Shape 103 (Time Series Analysis body textbox) : move + resize ---
$sh = $s.Shapes.Item(8)
$sh.Left = 723.153564453125
$sh.Top = 901.0039672851562
$sh.Width = 966.9448852539062
$sh.Height = 819.49609375

# --- Shape 105 (Linear Regression body textbox) : text fixes ---
$sh = $s.Shapes.Item(10)
$tr = $sh.TextFrame.TextRange
$tr.Paragraphs(6).Runs(1).Text = "Making transformations of predictor variables to find whether the response variable has a statistical interaction between the polynomial of quantitative variables and the interaction between qualitative variables and quantitative variables respectively as well as the interaction between quantitative variables themselves. After comparing each model’s AIC and Adjusted R-squared, we decide not to add any interaction terms in the model."
$tr.Paragraphs(8).Runs(1).Text = "The final model consists of the seven original predictor variables and the log transformation of response variable. The quantities of the final model can be concluded as:"

# --- Shape 107 (Conclusion body textbox) : text fix ---
$sh = $s.Shapes.Item(12)
$tr = $sh.TextFrame.TextRange
$tr.Paragraphs(3).Runs(1).Text = "In the linear regression, the model hasn’t been improved much by making any transformations of predictor variables or adding any interaction terms. This suggests that we may collect other variables like the automobile exhaust to see the deep link between PM2.5 and our life."

# --- Shape 108 (picture) : move + resize ---
$sh = $s.Shapes.Item(13)
$sh.Left = 1006.5
$sh.Top = 1624.23193359375
$sh.Width = 435.5256042480469
$sh.Height = 88.13024139404297

# --- Shape 109 (table) : change table style ---
$sh = $s.Shapes.Item(14)
$sh.Table.ApplyStyle("{A84C64D8-A4C9-479E-BBED-9062C8BFDFAE}")

# --- Shape 110 (picture) : move + resize ---
$sh = $s.Shapes.Item(15)
$sh.Left = 731.7106323242188
$sh.Top = 465.7027587890625
$sh.Width = 451.7892150878906
$sh.Height = 255.79718017578125

# --- Shape 111 (picture) : move + resize ---
$sh = $s.Shapes.Item(16)
$sh.Left = 1193.2579345703125
$sh.Top = 478.00396728515625
$sh.Width = 510.5157775878906
$sh.Height = 237.42323303222656

# --- Shape 112 ("Linear Regression" header) : move only ---
$sh = $s.Shapes.Item(17)
$sh.Left = 719.8651733398438
$sh.Top = 703.6624755859375
